$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "Bayesian Optimisation for XGBOOST"
$ws.Range("C10").Value = "Credit risk series"
$ws.Range("D10").Value = 4

$ws.Range("H15").Select()
